$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test data value in C2 from U7BM9E to G1S77A
$ws.Range("C2").Value = "G1S77A"

# Update the active cell selection to reflect where the edit was made
$ws.Range("D10").Select()
